$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.694.43"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.197.97"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.54"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.28"
$ws.Range("E6").Value = "  +5.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.199.69"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("E9").Value = "  +4.25%  "
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.97"
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.513"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.05"
$ws.Range("E14").Value = "  +4.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.723.48"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.677.44"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.43"
$ws.Range("E17").Value = "  +3.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.196.56"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.111"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "516.71"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.33"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.19"
$ws.Range("E23").Value = "  +4.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.97"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.18"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("E29").Value = "  +9.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.10"
$ws.Range("E30").Value = "  +8.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  +7.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.18"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.57"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "519.19"
$ws.Range("E36").Value = "  +9.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.96"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0899"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.126"
$ws.Range("E41").Value = "  +7.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.90"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0692"
$ws.Range("E43").Value = "  +13.41%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.305"
$ws.Range("E44").Value = "  +6.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.877.23"
$ws.Range("E46").Value = "  -4.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.68"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  +6.83%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.117"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.67"
$ws.Range("E51").Value = "  +9.43%  "
